$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    # Split the comma separated list of recorders and trim whitespace
    $rawParts = $text.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts.Count -le 1) {
        continue
    }

    # find the last index of an exact (case-sensitive) "System" entry
    $systemIdx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i] -ceq "System") {
            $systemIdx = $i
        }
    }

    if ($systemIdx -ge 0) {
        # move that "System" entry to the front, keep the rest in original order
        $rest = @()
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $systemIdx) {
                $rest += $parts[$i]
            }
        }
        $newParts = @("System") + $rest
    }
    else {
        # no exact "System" entry present - sort the names alphabetically
        $newParts = $parts | Sort-Object
    }

    $newText = [string]::Join(", ", $newParts)

    if ($newText -cne $text) {
        $ws.Cells.Item($r, 7).Value = $newText
    }
}
